$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade data rows to append (rows 9-11)
$data = @(
    @(9233.77, 9309.17, 107.89, 107.02, $false, -0.81, 42613.765497685185, $false),
    @(9196.83, 9233.77, 107.17, 106.74, $false, -0.4,  42614.672743055555, $false),
    @(9138.89, 9196.83, 107.04, 106.37, $false, -0.63, 42615.750092592592, $false)
)

$startRow = 9
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 7).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($r, 8).Value = $row[7]
}
